# Bewertung.xlsx: fill in group number + grade-self-assessment scores
# (D column) for the "Bewertung" sheet, matching the commit's grading pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bewertung")

# --- Data updates -----------------------------------------------------
# C3: Gruppenname (Nummer) -> 29
$ws.Range("C3").Value = 29

# D12: Funktionsumfang selbsteinschaetzung 8 -> 9
$ws.Range("D12").Value = 9

# D13: Verschluesselungsart 3 -> 5
$ws.Range("D13").Value = 5

# D14: Unittests + Abdeckung 9 -> (cleared)
$ws.Range("D14").ClearContents()

# D16: Type Checker 5 -> 0
$ws.Range("D16").Value = 0

# D18: Dokumentation 2 -> 5
$ws.Range("D18").Value = 5

# --- Selection change ---------------------------------------------------
$ws.Range("D7").Select()
